$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append two new data rows below the header ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 14
$ws.Range("C2").Value = 21307177
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "给了6分，1分是溢出来的"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 21307177
$ws.Range("D3").Value = 5

# --- Center-align the header + newly added data (matches new cellXfs style) ---
$ws.Range("A1:E2").HorizontalAlignment = -4108
$ws.Range("A3:D3").HorizontalAlignment = -4108

# --- Column widths (A:E) ---
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 13
$ws.Columns.Item(5).ColumnWidth = 25

# --- Selection moves to A4 after data entry ---
$ws.Range("A4").Select()
